# BoardGamesTemplate.xlsx — add a "Rok vydania" (Year of release) column
# between "Jednorázová" (I) and "Autori" (now shifted from J to K).
#
# Fixes bug reports #4, #5, #6:
#   - new template column + placeholder `${year}` for the release year
#   - AutoFilter / _FilterDatabase range grown to include the new column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at J; this pushes the existing "Autori" (J)
# and "Videonávod" (K) columns one slot to the right (K, L) and keeps
# their content/formatting intact.
$ws.Columns("J:J").Insert() | Out-Null

# Populate the new column's header + template placeholder row.
$ws.Range("J1").Value = "Rok vydania"
$ws.Range("J2").Value = '${year}'

# Match the look of the other narrow, centered columns (e.g. column A):
# width ~15 raw units (≈14.17 in COM "character" units) and centered text.
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Range("J1:J2").HorizontalAlignment = -4108 # xlCenter

# The AutoFilter used to cover C1:J2; grow it by one column to C1:K2 so it
# includes the newly inserted column.
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("C1:K2").AutoFilter() | Out-Null

# Re-point the workbook-level hidden _FilterDatabase name at the same
# widened range (toggling AutoFilter above recreates it, but still with
# the stale pre-insert bounds).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$C`$1:`$K`$2"
    }
}
